# Generate Report for Handback
#
# This script brings the workbook's localization-status report up to date
# after a handback: the "Status" text is updated everywhere it appears, the
# zh-cn / de-de detail sheets get their "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns filled in (with a
# hyperlink on the target-file cell, same as the source-file cell already
# has), and the columns that now hold longer text are widened to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/203a4c0b5c97a385479f4f9af815032102a510b4/e2e/"
$urlA = $baseUrl + "a.md"
$urlB = $baseUrl + "b.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells (columns E/F)
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Widen the status columns so the longer text fits (matches the width
# already used by other "fit" columns in this workbook).
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-16 12:35:47"
$zhcn.Range("K3").Value = "2016-08-16 12:35:47"

# Recreate the hyperlinks so the new "Latest Target File" link (column I)
# takes its place alongside the existing source-file link (column A).
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $urlA, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlA, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $urlB, "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlA, "", "", "a.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-08-16 12:35:55"
$dede.Range("K3").Value = "2016-08-16 12:35:55"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $urlA, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $urlA, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $urlB, "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $urlA, "", "", "a.md")

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Generated handback report"
